$wb = $excel.ActiveWorkbook

# Helper: force a value to be stored as literal text (avoids the engine's
# automatic ISO-date ("yyyy-mm-dd") parsing that Range.Value would trigger).
# We build a formula that evaluates to the literal text, then convert the
# cell to a static value via Copy + PasteSpecial(xlPasteValues), which does
# not disturb the cell's existing number format / style.
function Set-TextValue($range, [string]$text) {
    $range.Formula = '="' + $text + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

# ---------------------------------------------------------------------
# Sheet "AMSIN": append two new history rows (35 and 36) after the last
# existing row (34), copying the formatting of row 34 for the new rows.
# ---------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")

$wsAmsin.Range("A34:G34").Copy()
$wsAmsin.Range("A35:G36").PasteSpecial(-4122)  # xlPasteFormats

# Row 35
Set-TextValue $wsAmsin.Range("A35") "2023-02-17"
$wsAmsin.Range("B35").Value = 44974.45575663194
Set-TextValue $wsAmsin.Range("C35") "edu173fstccycle"
$wsAmsin.Range("D35").Value = 60
$wsAmsin.Range("E35").Value = 60
$wsAmsin.Range("F35").Value = 0
$wsAmsin.Range("G35").Value = 2.06

# Row 36
Set-TextValue $wsAmsin.Range("A36") "2023-02-20"
$wsAmsin.Range("B36").Value = 44977.42801012732
Set-TextValue $wsAmsin.Range("C36") "173educflow"
$wsAmsin.Range("D36").Value = 60
$wsAmsin.Range("E36").Value = 60
$wsAmsin.Range("F36").Value = 0
$wsAmsin.Range("G36").Value = 1.37

# ---------------------------------------------------------------------
# Sheet "AMS": row 22 picks up the same formatting as row 21, its B22
# run-time gets corrected, and two brand new rows (23, 24) are appended.
# ---------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")

# Row 22 adopts row 21's formatting (A/C/D/E/F/G become styled, B stays
# styled the same way it already was) and its run time is corrected.
$wsAms.Range("A21:G21").Copy()
$wsAms.Range("A22:G22").PasteSpecial(-4122)  # xlPasteFormats
$wsAms.Range("B22").Value = 44946.92204515047

# Row 23 (new) - fully styled like row 21/22.
$wsAms.Range("A21:G21").Copy()
$wsAms.Range("A23:G23").PasteSpecial(-4122)  # xlPasteFormats

Set-TextValue $wsAms.Range("A23") "2023-02-20"
$wsAms.Range("B23").Value = 44977.61088517361
Set-TextValue $wsAms.Range("C23") "173eduflow"
$wsAms.Range("D23").Value = 60
$wsAms.Range("E23").Value = 60
$wsAms.Range("F23").Value = 0
$wsAms.Range("G23").Value = 0.81

# Row 24 (new) - only column B carries the run-time style, like the
# original (pre-edit) row 22 did.
$wsAms.Range("B21").Copy()
$wsAms.Range("B24").PasteSpecial(-4122)  # xlPasteFormats

Set-TextValue $wsAms.Range("A24") "2023-02-20"
$wsAms.Range("B24").Value = 44977.82648602864
Set-TextValue $wsAms.Range("C24") "173educlive"
$wsAms.Range("D24").Value = 60
$wsAms.Range("E24").Value = 60
$wsAms.Range("F24").Value = 0
$wsAms.Range("G24").Value = 1.57
